# Applies the commit: adds a new "ODI Batting Extra" sheet (derived stats
# per match) and removes a handful of blank INNING_NUMBER placeholder cells
# on the "ODI Batting" sheet.

$wb = $excel.ActiveWorkbook

$wsPlayerInfo = $wb.Worksheets.Item("Player Info")
$wsOdiBatting = $wb.Worksheets.Item("ODI Batting")

# --- 1. Drop the empty INNING_NUMBER cells that used to be written out as
#        blank inline strings (rows where the player did not bat). ---
$wsOdiBatting.Range("B3").ClearContents()
$wsOdiBatting.Range("B4").ClearContents()
$wsOdiBatting.Range("B5").ClearContents()
$wsOdiBatting.Range("B9").ClearContents()

# --- 2. Add the new "ODI Batting Extra" worksheet after the existing
#        sheets so it lands as the third tab. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Header row - reuse the same bold/boxed header style used on the other
# sheets by copy/pasting the formats from an existing header cell.
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

$wsPlayerInfo.Range("A1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows. MATCH_CODE (A), NUM_4 (C), NUM_6 (D), PERCENT_RUNS_OF_TOTAL (E)
# and MAN_OF_MATCH (F) are textual; only BATTING_POSITION (B) is a genuine
# number, and it (plus NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL) is blank on rows
# where the player did not bat. A leading "'" forces Excel to keep a value
# as literal text instead of re-interpreting it as a number/percentage.

# Row 2 - match 3204
$newSheet.Range("A2").Value = "'3204"
$newSheet.Range("B2").Value = 7
$newSheet.Range("C2").Value = "'0"
$newSheet.Range("D2").Value = "'0"
$newSheet.Range("E2").Value = "'1.45%"
$newSheet.Range("F2").Value = "'NO"

# Row 3 - match 3206 (did not bat)
$newSheet.Range("A3").Value = "'3206"
$newSheet.Range("B3").Value = "'"
$newSheet.Range("C3").Value = "'"
$newSheet.Range("D3").Value = "'"
$newSheet.Range("E3").Value = "'"
$newSheet.Range("F3").Value = "'NO"

# Row 4 - match 3208 (did not bat)
$newSheet.Range("A4").Value = "'3208"
$newSheet.Range("B4").Value = 7
$newSheet.Range("C4").Value = "'"
$newSheet.Range("D4").Value = "'"
$newSheet.Range("E4").Value = "'"
$newSheet.Range("F4").Value = "'NO"

# Row 5 - match 3643 (did not bat)
$newSheet.Range("A5").Value = "'3643"
$newSheet.Range("B5").Value = 6
$newSheet.Range("C5").Value = "'"
$newSheet.Range("D5").Value = "'"
$newSheet.Range("E5").Value = "'"
$newSheet.Range("F5").Value = "'NO"

# Row 6 - match 3644
$newSheet.Range("A6").Value = "'3644"
$newSheet.Range("B6").Value = 6
$newSheet.Range("C6").Value = "'1"
$newSheet.Range("D6").Value = "'0"
$newSheet.Range("E6").Value = "'4.76%"
$newSheet.Range("F6").Value = "'NO"

# Row 7 - match 3645 (did not bat)
$newSheet.Range("A7").Value = "'3645"
$newSheet.Range("B7").Value = "'"
$newSheet.Range("C7").Value = "'"
$newSheet.Range("D7").Value = "'"
$newSheet.Range("E7").Value = "'"
$newSheet.Range("F7").Value = "'NO"

# Row 8 - match 3688
$newSheet.Range("A8").Value = "'3688"
$newSheet.Range("B8").Value = 6
$newSheet.Range("C8").Value = "'1"
$newSheet.Range("D8").Value = "'0"
$newSheet.Range("E8").Value = "'2.75%"
$newSheet.Range("F8").Value = "'NO"

# Row 9 - match 3689 (did not bat)
$newSheet.Range("A9").Value = "'3689"
$newSheet.Range("B9").Value = "'"
$newSheet.Range("C9").Value = "'"
$newSheet.Range("D9").Value = "'"
$newSheet.Range("E9").Value = "'"
$newSheet.Range("F9").Value = "'NO"

# Row 10 - match 3692 (did not bat)
$newSheet.Range("A10").Value = "'3692"
$newSheet.Range("B10").Value = "'"
$newSheet.Range("C10").Value = "'"
$newSheet.Range("D10").Value = "'"
$newSheet.Range("E10").Value = "'"
$newSheet.Range("F10").Value = "'NO"
